$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9148858785629272
$ws.Range("B1").Value = 1.580800175666809
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.484652280807495
$ws.Range("E1").Value = 1.382088184356689
